# "Add pie charts legend"
#
# The sheet held three side-by-side pie-chart data tables (Wojewodztwo /
# wytwarzanie / przesylanie) starting in column A. This edit makes room
# for a legend column: a new column A is inserted, the old column A
# (region names) becomes column B, and every other column shifts one to
# the right. The three header/category labels in the new/shifted column B
# are also renamed to the generic pie-chart-legend captions used by the
# chart, and column widths are adjusted to fit the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything one column to the right, freeing up column A.
$ws.Columns("A:A").Insert()

# Rename the (now shifted) header/category labels for the pie legend.
$ws.Range("B1").Value2 = "Nazwa wycinka"
$ws.Range("B2").Value2 = "Nakłady wiązane z wytwarzaniem"
$ws.Range("B3").Value2 = "Nakłady związane z przesyłem i dystrybucją"

# Resize columns to match the new layout: narrow new col A, wide label
# col B, and the data columns (C onward) slightly narrower than before.
$ws.Columns("A:A").ColumnWidth = 11.5204081632653
$ws.Columns("B:B").ColumnWidth = 36.5561224489796
$ws.Columns("C:AMK").ColumnWidth = 11.7704081632653

# Restore a sensible active selection on the sheet.
[void]$ws.Range("C8").Select()
